$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test data: the stored password value "Rajesh" becomes "Rajeshr"
$ws.Range("C2").Value = "Rajeshr"

# Update the phone-number style value in B2
$ws.Range("B2").Value = 8074453962

# Move the (non-data) active selection from E3 to N3
$ws.Range("N3").Select()
